# Casos de uso SUBE.xlsx - add new test cases (rows 21-24) for Tarjeta/Viaje
# use cases, and underline the header cell A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows 21-24 -----------------------------------------------
# Written in the exact order the strings first appear so the workbook's
# shared-string table is built up in the same sequence as the target file.

$ws.Range("D21").Value = "Devuelve si la tarjeta tiene o no la tarifa social"
$ws.Range("C22").Value = "tieneSaldoNegativo"
$ws.Range("D22").Value = "devuelve true o false si tiene saldo negativo"
$ws.Range("C23").Value = "cargaMinima"
$ws.Range("D23").Value = "devuelve el monto minimo que debe pagar para cargar la sube"
$ws.Range("C24").Value = "cerroViajeTren"
$ws.Range("D24").Value = "devuelve true o false si cerro el ultimo viaje en tren que hizo"
$ws.Range("C21").Value = "tieneTarifaSocial"

$ws.Range("B21").Value = "Tarjeta"
$ws.Range("F21").Value = "boolean"
$ws.Range("B22").Value = "Tarjeta"
$ws.Range("F22").Value = "boolean"
$ws.Range("B23").Value = "Tarjeta"
$ws.Range("F23").Value = "float"
$ws.Range("B24").Value = "Viaje"
$ws.Range("F24").Value = "boolean"

$ws.Range("A22").Value = 21
$ws.Range("A23").Value = 22
$ws.Range("A24").Value = 23

# --- Match formatting of the preceding rows (B, C, F columns) ----------
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B24").PasteSpecial(-4122)

$ws.Range("C20").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C24").PasteSpecial(-4122)

$ws.Range("F20").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("F24").PasteSpecial(-4122)

# --- Header cell A1 gets an underline -----------------------------------
$ws.Range("A1").Font.Underline = $true
$ws.Range("A1").HorizontalAlignment = 1

# --- Page setup (paper size / orientation) ------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Clean up the selection / scroll position left in the sheet view ---
$ws.Range("A1").Select()
